$wb = $excel.ActiveWorkbook

# --- Tasks sheet: add a "Cost" column (D) with value 1 for every task row ---
$wsTasks = $wb.Worksheets.Item("Tasks")
$wsTasks.Range("D1").Value = "Cost"
$wsTasks.Range("D2:D23").Value = 1

# --- Update selections / active sheet ---
# First touch the Resources sheet (it loses the "active tab" status but keeps
# an updated selection and no longer has a frozen/scrolled topLeftCell).
$wsResources = $wb.Worksheets.Item("Resources")
$wsResources.Activate()
$wsResources.Range("E33").Select()

# Finally activate Tasks and select D2 so it ends up as the active tab.
$wsTasks.Activate()
$wsTasks.Range("D2").Select()
